$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$deletedRow = 16

# Remove the "checksum" field row (the file-object property row that sat
# between "mime-type" and "file-size" under documents[].file). Deleting the
# entire row shifts everything below up by one; Excel auto-adjusts the
# merged "group name" cells in columns A and B to match (dimension becomes
# A1:N90).
#
# Known quirk: multi-row merges (e.g. "A44:A45") shift correctly, but
# single-CELL merges (e.g. "A42", a one-row group label) are left in place
# instead of moving up with their row. Capture the single-cell merges that
# sit below the deleted row beforehand, and re-home them afterwards so the
# final merge map lines up with the shifted data.
$singleCellMerges = New-Object System.Collections.ArrayList
foreach ($col in @("A", "B")) {
    for ($r = $deletedRow + 1; $r -le $ws.UsedRange.Rows.Count + 1; $r++) {
        $cell = $ws.Range("$col$r")
        if ($cell.MergeCells) {
            $area = $cell.MergeArea
            if ($area.Row -eq $r -and $area.Rows.Count -eq 1 -and $area.Columns.Count -eq 1) {
                [void]$singleCellMerges.Add("$col$r")
            }
        }
    }
}

$ws.Rows.Item($deletedRow).Delete()

foreach ($ref in $singleCellMerges) {
    if ($ref -match '^([A-Z]+)(\d+)$') {
        $col = $matches[1]
        $oldRow = [int]$matches[2]
        $newRow = $oldRow - 1
        $newRef = "$col$newRow"
        $target = $ws.Range($newRef)
        if (-not ($target.MergeCells -and $target.MergeArea.Address() -eq "`$$col`$$newRow")) {
            if ($target.MergeCells) {
                $target.MergeArea.UnMerge()
            }
            $target.Merge() | Out-Null
        }
    }
}
